$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting the old "RE_List" (D) column to E.
# This also pushes the dimension out to E14 and keeps the D1/E1 header style (s="1").
$ws.Columns.Item(4).Insert()

# New header text for the inserted column D1.
$ws.Range("D1").Value2 = "missingReplyNum"

# Re-run results: set the new missingReplyNum column (D2:D14) and refreshed RE_List (E2:E14).
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = "[0.0, 0.0, 0.0, 0.0, 1.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.25125628140703515, 9.022857573128268e-08, 10, 10, 0.0, 0.6153846153846154, 0.2222222222222222, 0.007400000000000018, 0.0, 0.625, 0.0, 0.020050125313283207, 0.002932551319648094, 0.0, 0.0, 0.0, 0.000408557836845193, 9.099182108918027, 0.0043098592890551774, 0.0, 0.010699999999999987, 9.140759373864446, 8.999920000639994, 0.0, 0.0, 0.08063730222209968, 0.9243842237319355, 0.9872093866965816, 2.141592653, 10, 10.0, 0.0015117157974301964, 1.7669115542410077, 0.04290000000000005, 1.0, 1.5, 0.4074, 1.0634202053261295, 10.0, 1.7146508254472612]"
$ws.Range("D3").Value2 = 0
$ws.Range("E3").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.1652802885849833, 3.6339005007031206e-05, 1.1953887114378692e-07, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0006518862940404747, 0.98997747998939, 0.010531821548960131, 0.019199999999999995, 0.0, 0.07043039072909245, 0.0, 0.0, 0.0, 0.1725735719998897, 0.9819981658328425, 0.9859483265353375, 2.141592653589793, 1.0, 0.8632, 0.005418304358780586, 10.0, 2.7071, 1.0, 1.5, 0.4226, 0.5415578112439925, 2.013502553319315, 0.10865510810250742]"
$ws.Range("D4").Value2 = 0
$ws.Range("E4").Value2 = "[0.0, 0.0, 0.0, 0.0, 1.3333333333333333, 0.3401360544217687, 0.0, 0.0, 0.0, 0.12224938875305623, 0.0, 0.0, 3.659517120546949e-11, 6.123917579087443e-08, 0.0, 0.0, 0.0, 0.0, 0.19799999999999995, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.9000399316735939, 8.82358738678542, 0.8993908923779054, 0.002439836221791475, 0.0, 0.01419192169581035, 5.479948160414716, 0.0, 0.0, 0.8387162019288229, 0.9739592533646291, 0.9584574469192972, 0.0, 1.0, 0.8632, 10.0, 0.9045892567503101, 1.2071, 1.0, 1.5, 0.4226, 0.43181428121307014, 0.6434696605587263, 0.06992470815776744]"
$ws.Range("D5").Value2 = 0
$ws.Range("E5").Value2 = "[0.0, 0.0, 0.0, 0.0, 10, 0.0, 0.0, 0.0, 0.0, 0.1295843520782396, 2.0, 0.9999183637573495, 10.0, 0.5075955023281672, 0.0, 0.0, 1.7777777777777777, 0.09259999999999996, 0.09899999999999998, 0.0, 0.0, 0.007518796992481203, 0.011730205278592375, 0.0, 0.0, 0.0, 0.9000412923264881, 9.22044238084454, 0.9897139073227031, 0.007200000000000012, 0.0, 0.01409847240727546, 0.0, 0.0, 0.0, 10, 0.9851252776098363, 0.9951406626680374, 2.0, 10, 10, 0.004495359127475493, 2.0531437839900772, 0.04290000000000005, 10.0, 1.2071, 0.5774, 0.13826690732881275, 0.05136677680985275, 0.6546245769151067]"
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0042796005706134095, 0.0, 0.02512562814070352, 6.441695647843136e-05, 2.129838963792373, 4.899134063392434, 0.0, 0.0, 0.0, 0.014800000000000035, 0.19799999999999995, 0.0, 0.0, 0.042606516290726815, 0.005865102639296188, 0.0, 0.0, 0.0, 0.00042600646578766515, 0.8986244154747128, 0.8998764716895982, 1.8072, 3.5455, 0.07043039072909245, 0.2800057599539204, 0.0, 0.0, 0.08063730222209968, 0.9247284013362598, 0.9581766315801489, 2.14, 1.0, 0.010199999999999987, 0.004787100025195156, 0.0018128041217441205, 0.15892540399838606, 1.0, 1.5, 0.8425999999999999, 0.34956176724067844, 0.5019525382997897, 0.03612626925467985]"
$ws.Range("D7").Value2 = 0
$ws.Range("E7").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.007132667617689016, 0.0, 0.25125628140703515, 0.00043458792466549554, 2.12982124626951, 0.5343423124865065, 1.2, 0.0, 0.0, 0.003999999999999948, 0.009000000000000008, 0.0, 0.0, 0.012531328320802004, 0.002932551319648094, 0.0, 0.0, 0.0, 0.9899819750124469, 0.9898819798594378, 0.9034317440515184, 0.0, 0.009500000000000008, 0.008459083444365221, 0.0, 0.0, 0.0, 0.39526160465565274, 10.0, 0.9951406626680374, 0.0, 1.4669624095260332, 0.01319999999999999, 0.0078105316200555485, 0.027573704799160357, 0.29290000000000005, 1.0, 1.5, 0.5774, 0.34956176724067844, 2.1045358966656655, 0.03612626925467985]"
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 3.6595171205469496e-06, 1.2247835158174887e-10, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.00026066097272970745, 0.0032388421147387238, 0.0038080339777700443, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.019637587224535912, 10.0, 10.0, 0.0, 0.0, 0.010599999999999943, 0.0, 0.21714531056196928, 0.20709999999999995, 1.0, 0.0, 0.4226, 0.0005853944094834694, 0.0007509762691498788, 0.0]"
$ws.Range("D9").Value2 = 0
$ws.Range("E9").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 2.251256281407035, 4.627106447758086e-11, 3.6595171205469496e-07, 9.859507302330785e-08, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.999000399840064, 0.899708219573986, 0.9006640965484581, 0.0, 0.0, 0.0, 0.0007919936640506004, 0.0, 0.0, 0.07564516277316556, 0.5096011662976776, 10.0, 2.14159, 0.0, 0.021094736842105277, 0.029982363315696564, 0.30434698592631143, 1.2071, 0.0, 0.366, 0.4074, 0.0, 0.0027035145689397236, 0.0]"
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 2.561661984382865e-07, 9.859507302330785e-08, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.00037039413980056074, 0.8997618014327166, 0.8999324799904715, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.03516562319000471, 10.0, 10.0, 0.0, 0.4669191437824923, 0.005199999999999982, 0.0, 0.39613586489838754, 1.4141, 0.0, 0.0, 0.4074, 0.0009268744816820395, 0.3496395313908081, 0.0015196518615735164]"
$ws.Range("D11").Value2 = 0
$ws.Range("E11").Value2 = "[1.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.734494844065125, 3.62292194934148e-06, 1.0222473254670867, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.03759398496240601, 0.0, 0.0, 0.0, 0.0, 0.9000390594602, 0.8990034048020872, 0.9003118000998371, 0.0015000000000000013, 0.0, 0.1126948222972764, 0.0, 0.0, 0.0, 0.053608038907429365, 10.0, 10.0, 0.0, 0.9338815533085252, 0.08420000000000005, 0.004661123708742663, 1.289857837992558, 1.7071, 0.0, 0.5, 0.4074, 0.04879912841276818, 0.3496395313908081, 0.32852110243835053]"
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.9000326050810861, 0.0002394825185443977, 0.001019554770536609, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.000983718086622347, 0.9851252776098363, 0.9951406626680374, 0.0, 0.46700443614289033, 0.0040999999999999925, 0.006487780297304232, 0.004865947905734201, 0.0, 0.0, 1.5, 0.0, 10.0, 0.00045058576148982056, 0.36077916695447954]"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 3.6595171205469496e-07, 9.79826812653991e-08, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 8.99640143942423, 0.0, 0.0, 9.8072, 0.0, 0.0, 2.239974080207358, 0.0, 0.0, 0.10434766619779166, 0.37939526302243487, 0.6305406791364982, 0.0, 0.46696240952603324, 0.0, 0.0, 0.0, 0.0, 0.0, 10, 0.0038000000000000256, 0.0006992211002162298, 0.0012015620306398329, 0.0005871382192443306]"
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = "[0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.16393442622950818, 0.0, 0.12224938875305623, 0.0, 0.0, 2.1300225196379494, 6.123916354303928, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.04887218045112782, 0.02932551319648094, 0.0, 0.0, 0.0, 0.9000390071273964, 0.899117747696819, 0.02907849159776472, 0.007200000000000012, 0.0, 0.0, 0.0, 0.0, 0.0, 0.6810087247520019, 0.5166322850745407, 0.6082789072442476, 0.0, 0.46700443614289033, 10.0, 0.004787100025195156, 0.008627857074706638, 0.20709999999999995, 0.0, 1.5, 0.4226, 0.1606094606240955, 0.08513667768098528, 0.35939766526213995]"
